# Update imputed KNN results on Sheet1 (Update Name of Algo)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -8.217000000000001
$ws.Range("B9").Value = 5.279999999999999
$ws.Range("D9").Value = -8.187999999999999
$ws.Range("D11").Value = -7.333
$ws.Range("B13").Value = 5.556
$ws.Range("B16").Value = 4.743
$ws.Range("D16").Value = -8.349
$ws.Range("B18").Value = 5.255999999999999
$ws.Range("B20").Value = 6.927
$ws.Range("D23").Value = -8.348000000000001
$ws.Range("D24").Value = -6.825
$ws.Range("B26").Value = 4.999
$ws.Range("D26").Value = -7.354000000000001
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("B29").Value = 5.677999999999999
$ws.Range("D34").Value = -7.536999999999999
$ws.Range("B35").Value = 9.35
$ws.Range("D35").Value = -8.013999999999999
$ws.Range("B36").Value = 8.079000000000001
$ws.Range("D44").Value = -7.323
$ws.Range("B45").Value = 5.689
$ws.Range("D48").Value = -7.56
$ws.Range("D49").Value = -8.35
$ws.Range("D52").Value = -7.923
$ws.Range("B55").Value = 4.813000000000001
$ws.Range("B57").Value = 5.232999999999999
$ws.Range("D66").Value = -7.316
$ws.Range("D67").Value = -7.485000000000001
$ws.Range("B69").Value = 5.712999999999999
$ws.Range("D73").Value = -8.042000000000002
$ws.Range("B76").Value = 6.545
$ws.Range("B78").Value = 8.409000000000002
$ws.Range("D78").Value = -8.141999999999999
$ws.Range("D80").Value = -7.994999999999999
$ws.Range("B82").Value = 5.456999999999999
$ws.Range("B83").Value = 5.529999999999999
$ws.Range("D91").Value = -7.498
$ws.Range("B93").Value = 5.831
$ws.Range("B97").Value = 4.999
$ws.Range("D97").Value = -7.345999999999999
$ws.Range("D99").Value = -8.253
$ws.Range("D104").Value = -7.907999999999999
